$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Se declara admisible, 05 días. Dese cuenta ONI."
$ws.Range("B2").Value = "Admisibles"

$ws.Range("A3").Value = "Pone en conc.inhabilidad"
$ws.Range("B3").Value = "Acumulación"

$ws.Range("A4").Value = "RESERVADO"
$ws.Range("B4").Value = "Ica Informa"

$ws.Range("A5").Value = "Dese cuenta admisibilidad"
$ws.Range("B5").Value = "Dese Cuenta"

$ws.Range("A6").Value = "Incompetencia"
$ws.Range("B6").Value = "Incompetencia"

$ws.Range("A7").Value = "Concede Recurso, Interconexión"
$ws.Range("B7").Value = "Concede Apelación"

$ws.Range("A8").Value = "Dese cuenta en sala la apelación"
$ws.Range("B8").Value = "Dese Cuenta"

$ws.Range("A9").Value = "Ev. Informe. En relación."
$ws.Range("B9").Value = "Evacua Informe"

$ws.Range("A10").Value = "Pide FUN dentro el palzo de 03 días"
$ws.Range("B10").Value = "Ica Solicita Diligencia"

$ws.Rows.Item(11).Delete()
